# Apply the commit's changes to the "Metadata" property sheet.
#
# Summary of the edit:
#   1. Update the "Date" property value (row 8, column B).
#   2. Insert a new "Jurisdiction" property row (with an empty value) right
#      after "Contact" (row 10) and before "Description" - this pushes
#      Description/Purpose/Copyright down by one row.
#   3. Append a new "Context" / "element:Element" property row at the end
#      of the table (row 21).
#
# The "Elements" sheet is left untouched - its cell contents do not change;
# only the shared-string table shifts as a side effect of the edits above
# (handled automatically by the engine when it serializes the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Bump the Date property to the new timestamp.
$ws.Range("B8").Value = "2024-09-12T14:01:50+00:00"

# 2) Insert a new blank row above "Description" (currently row 11) so it
#    becomes the new "Jurisdiction" row, then copy the formatting from the
#    row above ("Contact") so the new row matches the table's styling.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# 3) Append a new "Context" / "element:Element" row at the bottom of the
#    table (new row 21), copying the formatting of the last existing row.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Context"
$ws.Range("B21").Value = "element:Element"
